# "Worked on temporal resolution"
# The Demand sheet's single aggregate yearly value is replaced by a full
# 12-step time series (rows 3-14), and the Demand sheet becomes the
# active/selected sheet (previously SupIm was the active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# Row 3 (t=1): value changes from the old aggregate (3605375000) to the
# new per-step figure.
$ws.Cells.Item(3, 2).Value = 300447917

# Rows 4-14 (t=2..12): new rows with the same per-step demand value.
$demandValue = 300447917
for ($t = 2; $t -le 12; $t++) {
    $row = $t + 2
    $ws.Cells.Item($row, 1).Value = $t
    $ws.Cells.Item($row, 2).Value = $demandValue
}

# Selection moves to E14 on the Demand sheet, which also becomes the
# active/selected tab (it replaces SupIm as the active sheet).
[void]$ws.Range("E14").Select()
$ws.Activate()
